$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Constants_Details")
$ws.Range("A41").Value = "LOGIN_METADATA"
$ws.Range("B41").Value = "Login Metadata table."
$ws.Range("A42").Value = "TEACHER_METADATA"
$ws.Range("B42").Value = "Teacher metadata table."
$ws.Range("A43").Value = "STUDENT_METADATA"
$ws.Range("B43").Value = "Student Metadata table. "
